$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.424.09'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '2.524.29'
$ws.Range("E3").Value = '  -1.35%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.54'
$ws.Range("E5").Value = '  -2.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.95'
$ws.Range("E6").Value = '  +2.42%  '
$ws.Range("E7").Value = '  -1.32%  '
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.72'
$ws.Range("E10").Value = '  +0.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0802'
$ws.Range("E11").Value = '  -1.34%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.30'
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("E13").Value = '  -0.09%  '
$ws.Range("D14").Value = '2.913.98'
$ws.Range("E14").Value = '  -1.30%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.71'
$ws.Range("E15").Value = '  +4.11%  '
$ws.Range("D16").Value = '2.516.86'
$ws.Range("E16").Value = '  -4.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.818'
$ws.Range("E17").Value = '  -2.78%  '
$ws.Range("D18").Value = '42.409.51'
$ws.Range("E18").Value = '  -1.30%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.78'
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").Value = '0.0₃0952'
$ws.Range("E20").Value = '  -0.64%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.14'
$ws.Range("E21").Value = '  -3.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.98'
$ws.Range("E22").Value = '  -0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.34'
$ws.Range("E23").Value = '  -3.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.88'
$ws.Range("E24").Value = '  -2.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.03'
$ws.Range("E25").Value = '  -1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.85'
$ws.Range("E27").Value = '  -3.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.33'
$ws.Range("E28").Value = '  -4.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.17'
$ws.Range("E29").Value = '  -1.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.09'
$ws.Range("E30").Value = '  -0.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '155.75'
$ws.Range("E31").Value = '  +0.70%  '
$ws.Range("E32").Value = '  -1.62%  '
$ws.Range("E33").Value = '  +14.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0790'
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("E35").Value = '  -3.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.01'
$ws.Range("E36").Value = '  -5.08%  '
$ws.Range("B37").Value = 'Celestia'
$ws.Range("C37").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.13'
$ws.Range("E37").Value = '  -4.97%  '
$ws.Range("B38").Value = 'LidoDAOToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.16'
$ws.Range("E38").Value = '  -7.32%  '
$ws.Range("E39").Value = '  -0.68%  '
$ws.Range("E40").Value = '  +0.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.29'
$ws.Range("E41").Value = '  +8.93%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '21.83'
$ws.Range("E42").Value = '  -2.62%  '
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.29'
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("E45").Value = '  -2.38%  '
$ws.Range("D46").Value = '1.965.08'
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.88'
$ws.Range("E47").Value = '  -1.37%  '
$ws.Range("D48").Value = '2.769.48'
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '80.97'
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.854'
$ws.Range("E50").Value = '  +10.78%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.191'
$ws.Range("E51").Value = '  -0.76%  '
